$wb = $excel.ActiveWorkbook

# Set the final selection on "Azure Resources" first, so that after we
# finish working (and selecting) on "Most Important" below, that sheet
# remains the active/selected tab - matching the target workbook state.
$ws1 = $wb.Worksheets.Item("Azure Resources")
$ws1.Range("B132").Select()

# "Most Important" sheet: drop the old "Azure API Management" entry
# (row 1) and append two new entries at the bottom of the list.
$ws2 = $wb.Worksheets.Item("Most Important")
$ws2.Rows.Item(1).Delete()
$ws2.Range("A14").Value = "Azure Express Route"
$ws2.Range("A15").Value = "Azure API Manager"
$ws2.Range("A19").Select()
